$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 2
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 1

# Formula for I2 stays the same (H2/F2) - value recalculates automatically to 5
$ws.Range("I2").Formula = "=H2/F2"

# J2 formula changes from I2/G2 to I2/SQRT(G2)
$ws.Range("J2").Formula = "=I2/SQRT(G2)"

# Update the selected cell to H5
$ws.Range("H5").Select()
